# Applies the "Update to user data" commit: adds a 2016 data point for each
# Kenyan county alongside the existing 2013 figure.
#
# - Notes sheet: "Source: 45" -> "Source: 38"
# - Data sheet: for every existing (2013) county row, insert a new row right
#   below it holding the 2016 figure for that same county.
# - Data-wide-value sheet: add a "2016" column next to the existing "2013"
#   column, populated with the same per-county figures.

$wb = $excel.ActiveWorkbook

# 47 counties, in the same top-to-bottom order as the Data sheet, with their
# 2016 value.
$data = @(
    @("d18974", "Baringo", 48.4),
    @("d18975", "Bomet", 27.8),
    @("d18987", "Bungoma", 75.7),
    @("d18988", "Busia", 74.9),
    @("d18976", "Elgeyo-Marakwet", 44.2),
    @("d18955", "Embu", 65.9),
    @("d18965", "Garissa", 59.8),
    @("d18968", "Homa Bay", 33.9),
    @("d18956", "Isiolo", 75),
    @("d18991", "Kajiado", 87.3),
    @("d18989", "Kakamega", 89.9),
    @("d18977", "Kericho", 66.4),
    @("d18943", "Kiambu", 93.2),
    @("d18949", "Kilifi", 78.4),
    @("d18944", "Kirinyaga", 63.7),
    @("d18969", "Kisii", 89.4),
    @("d18970", "Kisumu", 79.9),
    @("d18957", "Kitui", 53.2),
    @("d18950", "Kwale", 60.6),
    @("d18978", "Laikipia", 64.2),
    @("d18951", "Lamu", 75),
    @("d18958", "Machakos", 68.4),
    @("d18959", "Makueni", 53.1),
    @("d18966", "Mandera", 32.8),
    @("d18960", "Marsabit", 42.9),
    @("d18961", "Meru", 70.2),
    @("d18971", "Migori", 60.9),
    @("d18952", "Mombasa", 82.6),
    @("d18946", "Murang'a", 67.8),
    @("d18964", "Nairobi", 97.1),
    @("d18979", "Nakuru", 72.3),
    @("d18980", "Nandi", 44.8),
    @("d18981", "Narok", 34.5),
    @("d18962", "Nithi", 59.3),
    @("d18972", "Nyamira", 83.6),
    @("d18947", "Nyandarua", 81.9),
    @("d18948", "Nyeri", 82.7),
    @("d18982", "Samburu", 42),
    @("d18973", "Siaya", 57.1),
    @("d18953", "Taita Taveta", 79.9),
    @("d18954", "Tana River", 67.2),
    @("d18983", "Trans-Nzoia", 79.8),
    @("d18984", "Turkana", 63.3),
    @("d18985", "Uasin Gishu", 69.3),
    @("d18990", "Vihiga", 88.3),
    @("d18967", "Wajir", 44.7),
    @("d18986", "West Pokot", 37.2)
)

# --- Notes sheet -----------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A4").Value = "Source: 38"

# --- Data sheet --------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("Data")

# Insert from the bottom up so row numbers above the insertion point never
# shift while we're still working on them.
for ($i = $data.Count - 1; $i -ge 0; $i--) {
    $row = $data[$i]
    $srcRow = $i + 2
    $newRow = $srcRow + 1

    $dataSheet.Rows.Item($newRow).Insert()

    $dataSheet.Cells.Item($newRow, 1).Value = $row[0]
    $dataSheet.Cells.Item($newRow, 2).Value = $row[1]
    $dataSheet.Cells.Item($newRow, 3).Value = 2016
    $dataSheet.Cells.Item($newRow, 4).Value = $row[2]
}

# --- Data-wide-value sheet --------------------------------------------------
$wideSheet = $wb.Worksheets.Item("Data-wide-value")

# A plain `.Value = "2016"` gets auto-coerced to the *number* 2016 by this
# engine (same heuristic real Excel applies to numeric-looking literals
# typed into a General-formatted cell). The existing "2013" header is a
# genuine text cell though, so match that: write it as a text formula, then
# collapse the formula to its literal (text) result with a values-only paste
# - this avoids touching the cell's number format/style entirely.
$wideSheet.Cells.Item(1, 3).Formula = "=""2016"""
$wideSheet.Cells.Item(1, 3).Copy()
$wideSheet.Cells.Item(1, 3).PasteSpecial(-4163)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $wideSheet.Cells.Item($r, 3).Value = $row[2]
}
